$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds rows of daily price data for
# "Femacal de La Calera - Repollo" from row 2 through row 848.
# Two new weekly records need to be inserted right before the
# existing row 732, which pushes all the following rows (732-848)
# down by two (to 734-850) and grows the used range from
# A1:R848 to A1:R850.

$ws.Rows("732:733").Insert()

# --- New row 732 ---
$ws.Cells.Item(732, 1).Value  = 3
$ws.Cells.Item(732, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(732, 3).Value  = "Coquimbo"
$ws.Cells.Item(732, 4).Value  = 44951
$ws.Cells.Item(732, 5).Value  = 5
$ws.Cells.Item(732, 6).Value  = 100112006
$ws.Cells.Item(732, 7).Value  = "Repollo"
$ws.Cells.Item(732, 8).Value  = "Crespo record"
$ws.Cells.Item(732, 9).Value  = "Primera"
$ws.Cells.Item(732, 10).Value = 2500
$ws.Cells.Item(732, 11).Value = 1100
$ws.Cells.Item(732, 12).Value = 1200
$ws.Cells.Item(732, 13).Value = 1152
$ws.Cells.Item(732, 14).Value = "`$/unidad"
$ws.Cells.Item(732, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(732, 16).Value = 1152
$ws.Cells.Item(732, 17).Value = 1
$ws.Cells.Item(732, 18).Value = "Hortaliza"

# --- New row 733 ---
$ws.Cells.Item(733, 1).Value  = 3
$ws.Cells.Item(733, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(733, 3).Value  = "Coquimbo"
$ws.Cells.Item(733, 4).Value  = 44951
$ws.Cells.Item(733, 5).Value  = 5
$ws.Cells.Item(733, 6).Value  = 100112006
$ws.Cells.Item(733, 7).Value  = "Repollo"
$ws.Cells.Item(733, 8).Value  = "Crespo record"
$ws.Cells.Item(733, 9).Value  = "Segunda"
$ws.Cells.Item(733, 10).Value = 1200
$ws.Cells.Item(733, 11).Value = 900
$ws.Cells.Item(733, 12).Value = 900
$ws.Cells.Item(733, 13).Value = 900
$ws.Cells.Item(733, 14).Value = "`$/unidad"
$ws.Cells.Item(733, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(733, 16).Value = 900
$ws.Cells.Item(733, 17).Value = 1
$ws.Cells.Item(733, 18).Value = "Hortaliza"
